$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Popularity Index for row 10 (I10): 510 -> 520
$ws.Range("I10").Value = 520

# Update ID for row 13 (C13): 207 -> 203
$ws.Range("C13").Value = 203

# Move the active cell selection to I16 (matches the saved cursor position)
$ws.Range("I16").Select()
